# Atualização de bases das ligas, do dia: 29-03-2024 às 13:24
#
# This script re-applies the corrected ordering for a handful of fixtures that
# were out of order (data rows 50/51, 89/90, 101/102/104 and 117/118), and
# appends the newly scraped fixture (row 125) for Banga Gargzdai vs
# FK Dziugas Telsiai.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Fix the ordering of a few rows: the odds/result data (columns B..AC)
#    belonged to the other row of each pair/group. Swap the payloads back in
#    place, keeping column A (the running id) fixed to the row position.
# ---------------------------------------------------------------------------

# Rows 50 <-> 51
$row50 = $ws.Range("B50:AC50").Value2
$row51 = $ws.Range("B51:AC51").Value2
$ws.Range("B50:AC50").Value2 = $row51
$ws.Range("B51:AC51").Value2 = $row50

# Rows 89 <-> 90
$row89 = $ws.Range("B89:AC89").Value2
$row90 = $ws.Range("B90:AC90").Value2
$ws.Range("B89:AC89").Value2 = $row90
$ws.Range("B90:AC90").Value2 = $row89

# Rows 101 -> 102 -> 104 -> 101 (3-way rotation)
$row101 = $ws.Range("B101:AC101").Value2
$row102 = $ws.Range("B102:AC102").Value2
$row104 = $ws.Range("B104:AC104").Value2
$ws.Range("B101:AC101").Value2 = $row104
$ws.Range("B102:AC102").Value2 = $row101
$ws.Range("B104:AC104").Value2 = $row102

# Rows 117 <-> 118
$row117 = $ws.Range("B117:AC117").Value2
$row118 = $ws.Range("B118:AC118").Value2
$ws.Range("B117:AC117").Value2 = $row118
$ws.Range("B118:AC118").Value2 = $row117

# ---------------------------------------------------------------------------
# 2) Append the new fixture as row 125 (Banga Gargzdai vs FK Dziugas Telsiai,
#    not played yet, so FTHG/FTAG/FTR and the PL_* result columns are blank).
# ---------------------------------------------------------------------------

# Copy the number formatting that the id (A) and Date (E) columns use from
# the last existing row so the new cells get the matching style.
$ws.Cells.Item(124, 1).Copy()
$ws.Cells.Item(125, 1).PasteSpecial(-4122)
$ws.Cells.Item(124, 5).Copy()
$ws.Cells.Item(125, 5).PasteSpecial(-4122)

$ws.Cells.Item(125, 1).Value2 = 123
$ws.Cells.Item(125, 2).Value2 = 7862915
$ws.Cells.Item(125, 3).Value2 = "Lithuania A Lyga"
$ws.Cells.Item(125, 4).Value2 = "Lithuania A Lyga"
$ws.Cells.Item(125, 5).Value2 = 45380.58333333334
$ws.Cells.Item(125, 6).Value2 = "Banga Gargzdai"
$ws.Cells.Item(125, 7).Value2 = "FK Dziugas Telsiai"
# H (FTHG), I (FTAG), J (FTR) intentionally left blank - match not played yet
$ws.Cells.Item(125, 11).Value2 = 2.6
$ws.Cells.Item(125, 12).Value2 = 2.9
$ws.Cells.Item(125, 13).Value2 = 2.625
$ws.Cells.Item(125, 14).Value2 = 2.7
$ws.Cells.Item(125, 15).Value2 = 2.8
$ws.Cells.Item(125, 16).Value2 = 2.7
$ws.Cells.Item(125, 17).Value2 = 0
$ws.Cells.Item(125, 18).Value2 = 1.9
$ws.Cells.Item(125, 19).Value2 = 1.9
$ws.Cells.Item(125, 20).Value2 = 2
$ws.Cells.Item(125, 21).Value2 = 1.95
$ws.Cells.Item(125, 22).Value2 = 1.85
$ws.Cells.Item(125, 23).Value2 = 0
$ws.Cells.Item(125, 24).Value2 = 0
$ws.Cells.Item(125, 25).Value2 = 0
$ws.Cells.Item(125, 26).Value2 = 0
$ws.Cells.Item(125, 27).Value2 = 0
# AB (PL_AhOver) and AC (PL_AhUnder) intentionally left blank
